$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.253.04"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "3.113.02"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'592.48"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'157.24"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "3.113.61"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("E10").Value = "  -5.70%  "
$ws.Range("D11").Value = "'5.92"
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'37.27"
$ws.Range("E13").Value = "  -4.90%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000241"
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "3.621.71"
$ws.Range("E16").Value = "  -23.85%  "
$ws.Range("D17").Value = "'7.25"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "64.071.31"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "3.107.36"
$ws.Range("E19").Value = "  -5.50%  "
$ws.Range("D20").Value = "'478.13"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  -7.46%  "
$ws.Range("D23").Value = "'7.60"
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("D24").Value = "'2.46"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("E25").Value = "  -5.25%  "
$ws.Range("D26").Value = "'81.18"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").Value = "'10.43"
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'7.50"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'2.69"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -5.78%  "
$ws.Range("D32").Value = "'2.20"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  -6.02%  "
$ws.Range("D34").Value = "'27.43"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").Value = "  -3.73%  "
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'6.04"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("E39").Value = "  -6.47%  "
$ws.Range("D40").Value = "'51.09"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").Value = "'9.24"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "'440.06"
$ws.Range("E42").Value = "  -8.95%  "
$ws.Range("D43").Value = "'0.293"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0365"
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").Value = "'40.40"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "2.830.16"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").Value = "'130.09"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "'25.49"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'2.25"
$ws.Range("E51").Value = "  -3.91%  "
